$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Epoch Accuracy values in column B (M05 Froze Decoder 3 results)
$ws.Range("B2").Value = 0.921875
$ws.Range("B4").Value = 0.8125
$ws.Range("B5").Value = 0.734375
$ws.Range("B6").Value = 0.6875
$ws.Range("B7").Value = 0.671875
$ws.Range("B8").Value = 0.65625
$ws.Range("B9").Value = 0.59375
$ws.Range("B10").Value = 0.640625
$ws.Range("B11").Value = 0.5625
$ws.Range("B12").Value = 0.59375
$ws.Range("B13").Value = 0.578125
$ws.Range("B14").Value = 0.671875
$ws.Range("B15").Value = 0.65625
$ws.Range("B16").Value = 0.578125
$ws.Range("B17").Value = 0.53125
$ws.Range("B18").Value = 0.515625
$ws.Range("B19").Value = 0.484375
$ws.Range("B20").Value = 0.5
$ws.Range("B21").Value = 0.515625
$ws.Range("B22").Value = 0.484375
$ws.Range("B23").Value = 0.515625
$ws.Range("B24").Value = 0.484375
$ws.Range("B25").Value = 0.484375
$ws.Range("B26").Value = 0.5
$ws.Range("B27").Value = 0.484375
$ws.Range("B28").Value = 0.484375
$ws.Range("B29").Value = 0.484375
$ws.Range("B30").Value = 0.46875
$ws.Range("B31").Value = 0.46875
$ws.Range("B32").Value = 0.46875
$ws.Range("B33").Value = 0.453125
$ws.Range("B34").Value = 0.46875
$ws.Range("B35").Value = 0.46875
$ws.Range("B36").Value = 0.453125
$ws.Range("B37").Value = 0.453125
$ws.Range("B38").Value = 0.453125
$ws.Range("B39").Value = 0.453125
$ws.Range("B40").Value = 0.453125
$ws.Range("B41").Value = 0.453125
$ws.Range("B42").Value = 0.453125
$ws.Range("B43").Value = 0.453125
$ws.Range("B44").Value = 0.453125
$ws.Range("B45").Value = 0.453125
$ws.Range("B46").Value = 0.453125
$ws.Range("B47").Value = 0.453125
$ws.Range("B48").Value = 0.453125
$ws.Range("B49").Value = 0.453125
$ws.Range("B50").Value = 0.453125
$ws.Range("B51").Value = 0.453125
$ws.Range("B52").Value = 0.453125
$ws.Range("B53").Value = 0.453125
$ws.Range("B54").Value = 0.453125
$ws.Range("B55").Value = 0.453125
$ws.Range("B56").Value = 0.453125
$ws.Range("B57").Value = 0.453125
$ws.Range("B58").Value = 0.453125
$ws.Range("B59").Value = 0.453125
$ws.Range("B60").Value = 0.453125
$ws.Range("B61").Value = 0.453125
$ws.Range("B62").Value = 0.453125
$ws.Range("B63").Value = 0.453125
$ws.Range("B64").Value = 0.453125
$ws.Range("B65").Value = 0.453125
$ws.Range("B66").Value = 0.453125
$ws.Range("B67").Value = 0.453125
$ws.Range("B68").Value = 0.453125
$ws.Range("B69").Value = 0.46875
$ws.Range("B70").Value = 0.46875
$ws.Range("B71").Value = 0.46875
$ws.Range("B72").Value = 0.46875
$ws.Range("B73").Value = 0.46875
$ws.Range("B74").Value = 0.46875
$ws.Range("B75").Value = 0.46875
$ws.Range("B76").Value = 0.46875
$ws.Range("B77").Value = 0.46875
$ws.Range("B78").Value = 0.46875
$ws.Range("B79").Value = 0.46875
$ws.Range("B80").Value = 0.484375
$ws.Range("B81").Value = 0.484375
$ws.Range("B82").Value = 0.484375
$ws.Range("B83").Value = 0.484375
$ws.Range("B84").Value = 0.484375
$ws.Range("B85").Value = 0.484375
$ws.Range("B86").Value = 0.484375
$ws.Range("B87").Value = 0.484375
$ws.Range("B88").Value = 0.484375
$ws.Range("B89").Value = 0.484375
$ws.Range("B90").Value = 0.484375
$ws.Range("B91").Value = 0.484375
$ws.Range("B92").Value = 0.484375
$ws.Range("B93").Value = 0.484375
$ws.Range("B94").Value = 0.484375
$ws.Range("B95").Value = 0.484375
$ws.Range("B96").Value = 0.484375
$ws.Range("B97").Value = 0.484375
$ws.Range("B98").Value = 0.484375
$ws.Range("B99").Value = 0.484375
$ws.Range("B100").Value = 0.484375
$ws.Range("B101").Value = 0.484375
$ws.Range("B102").Value = 0.484375
$ws.Range("B103").Value = 0.390625
$ws.Range("B104").Value = 0.5
$ws.Range("B105").Value = 0.40625
$ws.Range("B106").Value = 0.421875
$ws.Range("B107").Value = 0.46875
$ws.Range("B108").Value = 0.40625
$ws.Range("B109").Value = 0.515625
$ws.Range("B110").Value = 0.46875
$ws.Range("B111").Value = 0.375
$ws.Range("B112").Value = 0.328125
$ws.Range("B113").Value = 0.46875
$ws.Range("B114").Value = 0.46875
$ws.Range("B115").Value = 0.390625
$ws.Range("B116").Value = 0.453125
$ws.Range("B117").Value = 0.46875
$ws.Range("B118").Value = 0.4262295081967213

# Restore the active cell selection used when the workbook was last saved
$ws.Range("C9").Select()
